# Update automatico via Actualizar 03-07-2021 12-44-55
#
# This mirrors the "Actualizar disponibilidad" automation: a fresh check
# writes a new timestamp block at the top (rows 2-15) and every older
# timestamp block shifts down to the next block of 14 rows (16-29, then
# 30-43). Only column D (Fecha) changes; Nombre/URL/Disponibilidad stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D15").Value  = 44262.53097673065
$ws.Range("D16:D29").Value = 44262.50961724537
$ws.Range("D30:D43").Value = 44262.4882553588
